$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 107
$ws.Range("H107").Value = 2397.1428
$ws.Range("I107").Value = 2730
$ws.Range("J107").Value = 400
$ws.Range("K107").Value = 2730
$ws.Range("L107").Value = 400
$ws.Range("M107").Value = -810
$ws.Range("N107").Value = -4240

# Row 112
$ws.Range("H112").Value = 1913.75
$ws.Range("I112").Value = 700
$ws.Range("J112").Value = 2087.1428
$ws.Range("K112").Value = 2100
$ws.Range("L112").Value = 6261.428400000001
$ws.Range("M112").Value = -992
$ws.Range("N112").Value = -8477.428400000001

# Row 123
$ws.Range("H123").Value = 54980
$ws.Range("J123").Value = 54980
$ws.Range("L123").Value = 54980
$ws.Range("N123").Value = -64780

# Row 132
$ws.Range("H132").Value = 2115.95
$ws.Range("I132").Value = 2115.95
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 6347.849999999999
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -3817.849999999999
$ws.Range("N132").ClearContents()

# Row 137
$ws.Range("H137").Value = 52602864
$ws.Range("I137").Value = 12501128
$ws.Range("K137").Value = 37503384
$ws.Range("M137").Value = -37500834

$ws = $wb.Worksheets.Item("ARM")
# Row 81
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()

# Row 82
$ws.Range("H82").Value = 38400
$ws.Range("J82").Value = 38400
$ws.Range("L82").Value = 38400
$ws.Range("N82").Value = -39122

# Row 84
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()

# Row 85
$ws.Range("H85").Value = 38400
$ws.Range("J85").Value = 38400
$ws.Range("L85").Value = 38400
$ws.Range("N85").Value = -40896

# Row 130
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 1481.8064
$ws.Range("I94").Value = 1057.44
$ws.Range("J94").Value = 3250
$ws.Range("K94").Value = 1057.44
$ws.Range("L94").Value = 3250
$ws.Range("M94").Value = -606.4400000000001
$ws.Range("N94").Value = -4152

# Row 134
$ws.Range("H134").Value = 11905654
$ws.Range("I134").Value = 13158719
$ws.Range("J134").Value = 5103298.5
$ws.Range("K134").Value = 39476157
$ws.Range("L134").Value = 15309895.5
$ws.Range("M134").Value = -39473622
$ws.Range("N134").Value = -15314965.5

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 3160299
$ws.Range("I31").Value = 1668454.1
$ws.Range("J31").Value = 7822314
$ws.Range("K31").Value = 1668454.1
$ws.Range("L31").Value = 7822314
$ws.Range("M31").Value = -1668159.1
$ws.Range("N31").Value = -7822904

# Row 34
$ws.Range("H34").Value = 3160299
$ws.Range("I34").Value = 1668454.1
$ws.Range("J34").Value = 7822314
$ws.Range("K34").Value = 1668454.1
$ws.Range("L34").Value = 7822314
$ws.Range("M34").Value = -1668252.1
$ws.Range("N34").Value = -7822718

# Row 99
$ws.Range("H99").Value = 62514140
$ws.Range("I99").Value = 111120780
$ws.Range("J99").Value = 19878.572
$ws.Range("K99").Value = 111120780
$ws.Range("L99").Value = 19878.572
$ws.Range("M99").Value = -111119282
$ws.Range("N99").Value = -22874.572

# Row 126
$ws.Range("H126").Value = 62514140
$ws.Range("I126").Value = 111120780
$ws.Range("J126").Value = 19878.572
$ws.Range("K126").Value = 333362340
$ws.Range("L126").Value = 59635.716
$ws.Range("M126").Value = -333359870
$ws.Range("N126").Value = -64575.716

# Row 132
$ws.Range("H132").Value = 2002061.5
$ws.Range("I132").Value = 3126281
$ws.Range("J132").Value = 3448.889
$ws.Range("K132").Value = 9378843
$ws.Range("L132").Value = 10346.667
$ws.Range("M132").Value = -9376313
$ws.Range("N132").Value = -15406.667

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 3122157.5
$ws.Range("I5").Value = 2262923
$ws.Range("J5").Value = 5556655
$ws.Range("K5").Value = 6788769
$ws.Range("L5").Value = 16669965
$ws.Range("M5").Value = -6788657
$ws.Range("N5").Value = -16670189

# Row 47
$ws.Range("H47").Value = 262.5
$ws.Range("I47").Value = 262.5
$ws.Range("K47").Value = 787.5
$ws.Range("M47").Value = -356.5

# Row 55
$ws.Range("H55").Value = 871252.6
$ws.Range("I55").Value = 500
$ws.Range("J55").Value = 1054569
$ws.Range("K55").Value = 1500
$ws.Range("L55").Value = 3163707
$ws.Range("M55").Value = -1323
$ws.Range("N55").Value = -3164061

# Row 81
$ws.Range("H81").Value = 44028.62
$ws.Range("J81").Value = 47260.37
$ws.Range("L81").Value = 141781.11
$ws.Range("N81").Value = -144027.11

# Row 84
$ws.Range("H84").Value = 44028.62
$ws.Range("J84").Value = 47260.37
$ws.Range("L84").Value = 425343.33
$ws.Range("N84").Value = -436575.33

# Row 122
$ws.Range("H122").Value = 1561.5454
$ws.Range("I122").Value = 409.125
$ws.Range("K122").Value = 3682.125
$ws.Range("M122").Value = -1232.125

# Row 135
$ws.Range("H135").Value = 3122157.5
$ws.Range("I135").Value = 2262923
$ws.Range("J135").Value = 5556655
$ws.Range("K135").Value = 20366307
$ws.Range("L135").Value = 50009895
$ws.Range("M135").Value = -20363772
$ws.Range("N135").Value = -50014965

# Row 138
$ws.Range("H138").Value = 2950.3333
$ws.Range("I138").Value = 938.0909
$ws.Range("J138").Value = 5163.8
$ws.Range("K138").Value = 2814.2727
$ws.Range("L138").Value = 15491.4
$ws.Range("M138").Value = 2325.7273
$ws.Range("N138").Value = -25771.4

$ws = $wb.Worksheets.Item("GSM")
# Row 43
$ws.Range("H43").Value = 2400
$ws.Range("I43").Value = 1100
$ws.Range("J43").Value = 5000
$ws.Range("K43").Value = 1100
$ws.Range("L43").Value = 5000
$ws.Range("M43").Value = -949
$ws.Range("N43").Value = -5302

# Row 97
$ws.Range("H97").Value = 11906308
$ws.Range("I97").Value = 1137.9412
$ws.Range("J97").Value = 62503284
$ws.Range("K97").Value = 1137.9412
$ws.Range("L97").Value = 62503284
$ws.Range("M97").Value = -641.9412
$ws.Range("N97").Value = -62504276

$ws = $wb.Worksheets.Item("LTW")
# Row 74
$ws.Range("H74").Value = 24990
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 24990
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 24990
$ws.Range("M74").ClearContents()
$ws.Range("N74").Value = -26986

# Row 77
$ws.Range("H77").Value = 24990
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 24990
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 74970
$ws.Range("M77").ClearContents()
$ws.Range("N77").Value = -84954

$ws = $wb.Worksheets.Item("WVR")
# Row 34
$ws.Range("H34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("N34").ClearContents()

